# Split the run "de cuisine ou de table" + "</env></tl>" into four runs:
#   "de cuisine" | "</env>" | " ou de table" | "</tl>"
# preserving the exact character formatting of each piece (plain black text
# for the French words, Courier-New/blue for the XML-ish tag markup).

$d = $word.ActiveDocument

# --- Step 1: shrink "de cuisine ou de table" down to "de cuisine" -----------
# (this keeps the run's original, untouched rPr -- just less text in it)
$r1 = $d.Content
$r1.Find.Execute("de cuisine ou de table", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Text = "de cuisine"
$afterDeCuisine = $r1.End

# --- Step 2: locate the "</env>" markup that now immediately follows --------
# Bound the search to a short window right after our edit so we can't ever
# latch onto one of the other unrelated "</env>" tags elsewhere in the doc.
$r2 = $d.Range($afterDeCuisine, $afterDeCuisine + 20)
$r2.Find.Execute("</env>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# --- Step 3: pull "</env>" out (remembering its exact run formatting) -------
$r2.Cut()

# --- Step 4: insert the new plain-text segment " ou de table" ---------------
# Inserted directly after "de cuisine" it merges into that run, inheriting
# its clean black/no-rFonts formatting.
$insertPos = $d.Range($afterDeCuisine, $afterDeCuisine)
$insertPos.InsertBefore(" ou de table")

# --- Step 5: paste "</env>" back in, right before " ou de table" ------------
$pastePos = $d.Range($afterDeCuisine, $afterDeCuisine)
$pastePos.Paste()

# --- Step 6: force a clean re-write of the just-pasted run ------------------
# A straight Cut/Paste keeps the source run's old bookkeeping attributes
# (w:rsidDel/w:rsidR/w:rsidRPr); re-assigning its text (via a throwaway value
# first, so the assignment isn't a same-text no-op) makes the engine emit a
# fresh <w:r> the same way it would for genuinely new content, while leaving
# the run's rPr formatting untouched.
$pastedRange = $d.Range($afterDeCuisine, $afterDeCuisine + 6)
$pastedRange.Text = "zzzzzz"
$pastedRange2 = $d.Range($afterDeCuisine, $afterDeCuisine + 6)
$pastedRange2.Text = "</env>"
